$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.164.29'
$ws.Range("D3").Value = '1.829.23'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("D4").Value = "'0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = "'242.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").Value = "'0.6201"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.84%  '
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").Value = "'0.07347"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.72%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = "'0.2900"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.86%  '
$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").Value = "'23.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.72%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = "'0.07632"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.836.74'
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'4.962"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.69%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = "'0.6688"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.45%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = "'82.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.32%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = "'0.000008975"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.84%  '
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").Value = "'5.834"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.28%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '29.151.98'
$ws.Range("E18").Value = '  +0.40%  '
$ws.Range("B19").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C19").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D19").Value = '2.084.77'
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = "'235.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.59%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").Value = "'12.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.56%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = "'1.0000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").Value = "'7.345"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.12%  '
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").Value = "'1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = "'158.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = "'0.1390"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.47%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = "'8.515"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.36%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = "'17.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.11%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = "'1.485"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.90%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = "'0.05847"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.31%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = "'4.076"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'4.083"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.64%  '
$ws.Range("B33").Value = 'Toncoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D33").Value = "'1.212"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.79%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = "'1.861"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.58%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = "'1.137"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = "'0.7247"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.88%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = "'2.608"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.00%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = "'2.844"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.55%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.226.55'
$ws.Range("E39").Value = '  +0.91%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = "'0.01758"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.14%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = "'6.221"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.21%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = "'0.9078"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.03%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = "'1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = "'101.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").Value = '1.985.87'
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'65.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.39%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = "'0.5036"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.03%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = "'9.142"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.04%  '
$ws.Range("D49").Value = "'0.4032"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.87%  '
$ws.Range("E50").Value = '  -2.21%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = "'0.1129"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.87%  '
